$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply header style (same as existing header cells, e.g. H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$data = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(8, 8)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(10, 11)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(8, 8)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(6, 6)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(6, 6)
    29 = @(5, 5)
    30 = @(9, 9)
    31 = @(8, 8)
    32 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
